$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) keeps its original text representation instead of
# being auto-converted to a floating point number by Excel (values like
# "578.48" or "0.0000276" must remain exact text, matching the source data).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '67.152.96'
$ws.Range('E2').Value = '  +2.57%  '
$ws.Range('D3').Value = '3.445.00'
$ws.Range('E3').Value = '  +1.61%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '578.48'
$ws.Range('E5').Value = '  +3.41%  '
$ws.Range('D6').Value = '187.36'
$ws.Range('E6').Value = '  +6.61%  '
$ws.Range('D7').Value = '0.629'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '3.440.32'
$ws.Range('E8').Value = '  +1.70%  '
$ws.Range('D9').Value = '0.999'
$ws.Range('E9').Value = '  -0.10%  '
$ws.Range('D10').Value = '0.172'
$ws.Range('E10').Value = '  +0.71%  '
$ws.Range('D11').Value = '0.641'
$ws.Range('E11').Value = '  +0.57%  '
$ws.Range('D12').Value = '57.79'
$ws.Range('E12').Value = '  +8.18%  '
$ws.Range('D13').Value = '0.0000276'
$ws.Range('E13').Value = '  -0.77%  '
$ws.Range('E14').Value = '  +2.34%  '
$ws.Range('D15').Value = '3.993.74'
$ws.Range('E15').Value = '  +1.52%  '
$ws.Range('D16').Value = '18.95'
$ws.Range('E16').Value = '  +3.57%  '
$ws.Range('D17').Value = '3.434.12'
$ws.Range('E17').Value = '  +1.05%  '
$ws.Range('D18').Value = '66.947.22'
$ws.Range('E18').Value = '  +2.16%  '
$ws.Range('E19').Value = '  -0.58%  '
$ws.Range('D20').Value = '12.04'
$ws.Range('E20').Value = '  +1.76%  '
$ws.Range('D21').Value = '1.02'
$ws.Range('E21').Value = '  +1.63%  '
$ws.Range('D22').Value = '490.38'
$ws.Range('E22').Value = '  +5.26%  '
$ws.Range('D23').Value = '5.60'
$ws.Range('E23').Value = '  +13.40%  '
$ws.Range('D24').Value = '16.94'
$ws.Range('E24').Value = '  +18.41%  '
$ws.Range('D25').Value = '4.32'
$ws.Range('E25').Value = '  +4.77%  '
$ws.Range('D26').Value = '89.44'
$ws.Range('E26').Value = '  +2.53%  '
$ws.Range('D27').Value = '2.96'
$ws.Range('E27').Value = '  +1.62%  '
$ws.Range('D28').Value = '10.92'
$ws.Range('E28').Value = '  +1.95%  '
$ws.Range('D29').Value = '9.00'
$ws.Range('E29').Value = '  +3.14%  '
$ws.Range('D30').Value = '31.18'
$ws.Range('E30').Value = '  +0.33%  '
$ws.Range('D31').Value = '7.32'
$ws.Range('E31').Value = '  +11.88%  '
$ws.Range('D32').Value = '604.56'
$ws.Range('E32').Value = '  +5.37%  '
$ws.Range('D33').Value = '64.82'
$ws.Range('E33').Value = '  +2.46%  '
$ws.Range('D34').Value = '11.78'
$ws.Range('E34').Value = '  +2.58%  '
$ws.Range('D35').Value = '0.112'
$ws.Range('E35').Value = '  +3.42%  '
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('D37').Value = '0.146'
$ws.Range('E37').Value = '  +2.86%  '
$ws.Range('D38').Value = '36.92'
$ws.Range('E38').Value = '  +2.80%  '
$ws.Range('D39').Value = '0.0₃0783'
$ws.Range('E39').Value = '  +5.79%  '
$ws.Range('D40').Value = '0.386'
$ws.Range('E40').Value = '  +3.39%  '
$ws.Range('D41').Value = '3.45'
$ws.Range('E41').Value = '  -4.50%  '
$ws.Range('D42').Value = '3.190.39'
$ws.Range('E42').Value = '  +2.18%  '
$ws.Range('D43').Value = '2.88'
$ws.Range('E43').Value = '  +3.10%  '
$ws.Range('D44').Value = '0.0429'
$ws.Range('E44').Value = '  +2.89%  '
$ws.Range('D45').Value = '2.56'
$ws.Range('E45').Value = '  +5.05%  '
$ws.Range('D46').Value = '3.23'
$ws.Range('E46').Value = '  +1.48%  '
$ws.Range('E47').Value = '  +1.14%  '
$ws.Range('D48').Value = '2.65'
$ws.Range('E48').Value = '  +14.99%  '
$ws.Range('D49').Value = '0.997'
$ws.Range('E49').Value = '  -0.23%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').Value = '8.64'
$ws.Range('E50').Value = '  +2.49%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').Value = '140.03'
$ws.Range('E51').Value = '  -0.06%  '
